$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to text format so numeric-looking values stay as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "67.227.31"
$ws.Cells.Item(2, 5).Value = "  -8.36%  "
$ws.Cells.Item(3, 4).Value = "3.667.11"
$ws.Cells.Item(3, 5).Value = "  -7.85%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "562.88"
$ws.Cells.Item(5, 5).Value = "  -7.80%  "
$ws.Cells.Item(6, 4).Value = "168.76"
$ws.Cells.Item(6, 5).Value = "  +0.02%  "
$ws.Cells.Item(7, 4).Value = "3.656.52"
$ws.Cells.Item(7, 5).Value = "  -7.91%  "
$ws.Cells.Item(8, 4).Value = "0.613"
$ws.Cells.Item(8, 5).Value = "  -10.19%  "
$ws.Cells.Item(9, 4).Value = "0.999"
$ws.Cells.Item(9, 5).Value = "  -0.08%  "
$ws.Cells.Item(10, 4).Value = "0.691"
$ws.Cells.Item(10, 5).Value = "  -12.43%  "
$ws.Cells.Item(11, 4).Value = "0.160"
$ws.Cells.Item(11, 5).Value = "  -14.51%  "
$ws.Cells.Item(12, 4).Value = "50.68"
$ws.Cells.Item(12, 5).Value = "  -9.99%  "
$ws.Cells.Item(13, 4).Value = "0.0000287"
$ws.Cells.Item(13, 5).Value = "  -14.79%  "
$ws.Cells.Item(14, 4).Value = "10.31"
$ws.Cells.Item(14, 5).Value = "  -8.36%  "
$ws.Cells.Item(15, 4).Value = "4.278.74"
$ws.Cells.Item(15, 5).Value = "  -7.32%  "
$ws.Cells.Item(16, 4).Value = "3.716.55"
$ws.Cells.Item(16, 5).Value = "  -6.58%  "
$ws.Cells.Item(17, 5).Value = "  -3.42%  "
$ws.Cells.Item(18, 4).Value = "18.99"
$ws.Cells.Item(18, 5).Value = "  -8.91%  "
$ws.Cells.Item(19, 5).Value = "  -10.55%  "
$ws.Cells.Item(20, 4).Value = "12.60"
$ws.Cells.Item(20, 5).Value = "  -11.94%  "
$ws.Cells.Item(21, 4).Value = "67.070.50"
$ws.Cells.Item(21, 5).Value = "  -8.52%  "
$ws.Cells.Item(22, 4).Value = "398.66"
$ws.Cells.Item(22, 5).Value = "  -12.11%  "
$ws.Cells.Item(23, 4).Value = "4.40"
$ws.Cells.Item(23, 5).Value = "  -8.89%  "
$ws.Cells.Item(24, 4).Value = "86.41"
$ws.Cells.Item(24, 5).Value = "  -10.08%  "
$ws.Cells.Item(25, 4).Value = "2.99"
$ws.Cells.Item(25, 5).Value = "  -12.16%  "
$ws.Cells.Item(26, 4).Value = "12.43"
$ws.Cells.Item(26, 5).Value = "  -12.66%  "
$ws.Cells.Item(27, 4).Value = "10.37"
$ws.Cells.Item(27, 5).Value = "  -6.32%  "
$ws.Cells.Item(28, 5).Value = "  -0.29%  "
$ws.Cells.Item(29, 4).Value = "3.63"
$ws.Cells.Item(29, 5).Value = "  -13.08%  "
$ws.Cells.Item(30, 4).Value = "9.27"
$ws.Cells.Item(30, 5).Value = "  -11.97%  "
$ws.Cells.Item(31, 4).Value = "32.07"
$ws.Cells.Item(31, 5).Value = "  -11.85%  "
$ws.Cells.Item(32, 4).Value = "7.41"
$ws.Cells.Item(32, 5).Value = "  -7.45%  "
$ws.Cells.Item(33, 4).Value = "12.29"
$ws.Cells.Item(33, 5).Value = "  -12.01%  "
$ws.Cells.Item(34, 4).Value = "0.114"
$ws.Cells.Item(34, 5).Value = "  -12.06%  "
$ws.Cells.Item(35, 4).Value = "63.86"
$ws.Cells.Item(35, 5).Value = "  -9.75%  "
$ws.Cells.Item(36, 4).Value = "42.20"
$ws.Cells.Item(36, 5).Value = "  -11.96%  "
$ws.Cells.Item(37, 4).Value = "575.36"
$ws.Cells.Item(37, 5).Value = "  -11.36%  "
$ws.Cells.Item(38, 4).Value = "0.0₃0868"
$ws.Cells.Item(38, 5).Value = "  -18.17%  "
$ws.Cells.Item(39, 4).Value = "0.999"
$ws.Cells.Item(39, 5).Value = "  -0.10%  "
$ws.Cells.Item(40, 5).Value = "  +0.17%  "
$ws.Cells.Item(41, 4).Value = "0.387"
$ws.Cells.Item(41, 5).Value = "  -10.29%  "
$ws.Cells.Item(42, 4).Value = "0.132"
$ws.Cells.Item(42, 5).Value = "  -10.22%  "
$ws.Cells.Item(43, 4).Value = "2.95"
$ws.Cells.Item(43, 5).Value = "  -12.64%  "
$ws.Cells.Item(44, 4).Value = "0.0428"
$ws.Cells.Item(44, 5).Value = "  -11.53%  "
$ws.Cells.Item(45, 5).Value = "  -13.73%  "
$ws.Cells.Item(46, 4).Value = "2.50"
$ws.Cells.Item(46, 5).Value = "  -3.61%  "
$ws.Cells.Item(47, 4).Value = "9.00"
$ws.Cells.Item(47, 5).Value = "  -15.46%  "
$ws.Cells.Item(48, 4).Value = "0.132"
$ws.Cells.Item(48, 5).Value = "  -11.39%  "
$ws.Cells.Item(49, 4).Value = "2.66"
$ws.Cells.Item(49, 5).Value = "  -15.26%  "
# Rows 50 and 51 swap content: Maker <-> ApeXProtocol plus updated values
$ws.Cells.Item(50, 2).Value = "ApeXProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(50, 4).Value = "3.12"
$ws.Cells.Item(50, 5).Value = "  -10.19%  "
$ws.Cells.Item(51, 2).Value = "Maker"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(51, 4).Value = "2.693.47"
$ws.Cells.Item(51, 5).Value = "  -4.32%  "
